$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new data row at row 436, pushing the existing rows
# 436-478 down to 437-479 (this also grows the used range from R478 to R479).
$ws.Rows.Item(436).Insert()

# Populate the newly inserted row 436 with the new weekly record.
# All "dimension" columns (market, region, category, etc.) are constant
# across this sheet's series; only the date/price columns vary per row.
$ws.Cells.Item(436, 1).Value = 8
$ws.Cells.Item(436, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(436, 3).Value = "Coquimbo"
$ws.Cells.Item(436, 4).Value = 45194
$ws.Cells.Item(436, 5).Value = 4
$ws.Cells.Item(436, 6).Value = 100112012
$ws.Cells.Item(436, 7).Value = "Espinaca"
$ws.Cells.Item(436, 8).Value = "Sin especificar"
$ws.Cells.Item(436, 9).Value = "Primera"
$ws.Cells.Item(436, 10).Value = 1400
$ws.Cells.Item(436, 11).Value = 450
$ws.Cells.Item(436, 12).Value = 500
$ws.Cells.Item(436, 13).Value = 475
$ws.Cells.Item(436, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(436, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(436, 16).Value = 950
$ws.Cells.Item(436, 17).Value = 0.5
$ws.Cells.Item(436, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date number format used by the
# rest of column D in this sheet.
$ws.Cells.Item(436, 4).NumberFormat = $ws.Cells.Item(437, 4).NumberFormat
